$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Corrections to existing CDC cumulative case counts
$ws.Range("B43").Value = 78
$ws.Range("B46").Value = 211
$ws.Range("B47").Value = 275
$ws.Range("B48").Value = 422
$ws.Range("B55").Value = 3471
$ws.Range("B57").Value = 7023
$ws.Range("B63").Value = 44338
$ws.Range("B74").Value = 258098
$ws.Range("B75").Value = 267436
$ws.Range("B78").Value = 395926

# Extend column A's date formatting down to the new rows (matches existing style)
$ws.Range("A94").Copy()
$ws.Range("A95:A97").PasteSpecial(-4122)

# Append new rows for the latest reporting dates
$ws.Range("A95").Value = 43945
$ws.Range("B95").Value = 895766

$ws.Range("A96").Value = 43946
$ws.Range("B96").Value = 928619

$ws.Range("A97").Value = 43947
$ws.Range("B97").Value = 957875
